$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: 23/05/2018 - Borghesan Simona - Silesia Nera - Mt. - 1
$ws.Range("A32").Value = "05/23/2018"
$ws.Range("A32").NumberFormat = "dd/mm/yyyy"
$ws.Range("A32").VerticalAlignment = -4160

$ws.Range("B32").Value = "Borghesan Simona"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").VerticalAlignment = -4160

$ws.Range("C32").Value = "Silesia Nera"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").VerticalAlignment = -4160

$ws.Range("D32").Value = "Mt."
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").VerticalAlignment = -4160

$ws.Range("E32").Value = 1
$ws.Range("E32").VerticalAlignment = -4160

# Row 33: 24/05/2018 - Licata Rosa - Rocche filo nero - N°. - 5
$ws.Range("A33").Value = "05/24/2018"
$ws.Range("A33").NumberFormat = "dd/mm/yyyy"
$ws.Range("A33").VerticalAlignment = -4160

$ws.Range("B33").Value = "Licata Rosa"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").VerticalAlignment = -4160

$ws.Range("C33").Value = "Rocche filo nero"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").VerticalAlignment = -4160

$ws.Range("D33").Value = "N°."
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").VerticalAlignment = -4160

$ws.Range("E33").Value = 5
$ws.Range("E33").VerticalAlignment = -4160
